# DDAf_2025_Tableau_annexe_Tab26.xlsx -- "Add files via upload" re-edit
# Refreshes the Table 26 (Subjective well-being) data block (rows 4-98, cols C:F)
# with the latest World Happiness Report 2025 figures, and fixes a mojibake
# (UTF-8-decoded-as-Latin-1) footnote string about the PALOP regional grouping.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C4").Value = 3.438
$ws.Range("D4").Value = 0.697
$ws.Range("E4").Value = 0.018
$ws.Range("F4").Value = 1.146

$ws.Range("C5").Value = 3.774
$ws.Range("D5").Value = 0.428
$ws.Range("E5").Value = 0.054
$ws.Range("F5").Value = 1.152

$ws.Range("C6").Value = 3.757
$ws.Range("D6").Value = 0.606
$ws.Range("E6").Value = 0.065
$ws.Range("F6").Value = 1.131

$ws.Range("C7").Value = 3.26
$ws.Range("D7").Value = 0.677
$ws.Range("E7").Value = 0.115
$ws.Range("F7").Value = 0.489

$ws.Range("C8").Value = 5.19
$ws.Range("D8").Value = 0.813
$ws.Range("E8").Value = 0.127
$ws.Range("F8").Value = 1.131

$ws.Range("C9").Value = 4.911
$ws.Range("D9").Value = 0.564
$ws.Range("E9").Value = 0.052
$ws.Range("F9").Value = 1.482

$ws.Range("C10").Value = 5.213
$ws.Range("D10").Value = 0.676
$ws.Range("E10").Value = 0.064
$ws.Range("F10").Value = 1.465

$ws.Range("C11").Value = 3.912
$ws.Range("D11").Value = 0.872
$ws.Range("E11").Value = 0.131
$ws.Range("F11").Value = 1.013

$ws.Range("C12").Value = 3.396
$ws.Range("D12").Value = 0.598
$ws.Range("E12").Value = 0.065
$ws.Range("F12").Value = 0.961

$ws.Range("C13").Value = 4.0646
$ws.Range("D13").Value = 0.659
$ws.Range("E13").Value = 0.076778
$ws.Range("F13").Value = 1.107778

$ws.Range("C15").Value = 4.887
$ws.Range("D15").Value = 0.662
$ws.Range("E15").Value = 0.099
$ws.Range("F15").Value = 0.986

$ws.Range("C17").Value = 4.384
$ws.Range("D17").Value = 0.477
$ws.Range("E17").Value = 0.175
$ws.Range("F17").Value = 0.902

$ws.Range("C18").Value = 5.03
$ws.Range("D18").Value = 0.626
$ws.Range("E18").Value = 0.082
$ws.Range("F18").Value = 0.796

$ws.Range("C19").Value = 3.469
$ws.Range("D19").Value = 0.6
$ws.Range("E19").Value = 0.151
$ws.Range("F19").Value = 0.929

$ws.Range("C21").Value = 5.12
$ws.Range("D21").Value = 0.651
$ws.Range("E21").Value = 0.042
$ws.Range("F21").Value = 1.224

$ws.Range("C23").Value = 4.305843
$ws.Range("D23").Value = 0.515259
$ws.Range("E23").Value = 0.139216
$ws.Range("F23").Value = 0.699857

$ws.Range("C24").Value = 3.754
$ws.Range("D24").Value = 0.26
$ws.Range("E24").Value = 0.117
$ws.Range("F24").Value = 0.501

$ws.Range("C27").Value = 3.898
$ws.Range("D27").Value = 0.54
$ws.Range("E27").Value = 0.182
$ws.Range("F27").Value = 1.051

$ws.Range("C28").Value = 4.51
$ws.Range("D28").Value = 0.649
$ws.Range("E28").Value = 0.22
$ws.Range("F28").Value = 1.047

$ws.Range("C29").Value = 4.157
$ws.Range("D29").Value = 0.292
$ws.Range("E29").Value = 0.131
$ws.Range("F29").Value = 0.938

$ws.Range("C30").Value = 5.832
$ws.Range("D30").Value = 0.761
$ws.Range("E30").Value = 0.095
$ws.Range("F30").Value = 1.521

$ws.Range("C33").Value = 4.347
$ws.Range("D33").Value = 0.916
$ws.Range("E33").Value = 0.206
$ws.Range("F33").Value = 0.877

$ws.Range("C36").Value = 3.8
$ws.Range("D36").Value = 0.858
$ws.Range("E36").Value = 0.153
$ws.Range("F36").Value = 0.908

$ws.Range("C37").Value = 4.461
$ws.Range("D37").Value = 0.751
$ws.Range("E37").Value = 0.148
$ws.Range("F37").Value = 1.342

$ws.Range("C38").Value = 4.089327
$ws.Range("D38").Value = 0.571361
$ws.Range("E38").Value = 0.164893
$ws.Range("F38").Value = 0.887128

$ws.Range("C39").Value = 5.571
$ws.Range("D39").Value = 0.522
$ws.Range("E39").Value = 0.057
$ws.Range("F39").Value = 1.363

$ws.Range("C40").Value = 3.817
$ws.Range("D40").Value = 0.593
$ws.Range("E40").Value = 0.018
$ws.Range("F40").Value = 1.075

$ws.Range("C41").Value = 5.82
$ws.Range("D41").Value = 0.724
$ws.Range("E41").Value = 0.111
$ws.Range("F41").Value = 1.309

$ws.Range("C42").Value = 4.542
$ws.Range("D42").Value = 0.454
$ws.Range("E42").Value = 0.113
$ws.Range("F42").Value = 0.919

$ws.Range("C43").Value = 4.622
$ws.Range("D43").Value = 0.748
$ws.Range("E43").Value = 0.031
$ws.Range("F43").Value = 0.635

$ws.Range("C44").Value = 4.552
$ws.Range("D44").Value = 0.382
$ws.Range("E44").Value = 0.032
$ws.Range("F44").Value = 1.224

$ws.Range("C45").Value = 4.820667
$ws.Range("D45").Value = 0.5705
$ws.Range("E45").Value = 0.060333
$ws.Range("F45").Value = 1.0875

$ws.Range("C46").Value = 4.357
$ws.Range("D46").Value = 0.679
$ws.Range("E46").Value = 0.092
$ws.Range("F46").Value = 0.228

$ws.Range("C47").Value = 4.383
$ws.Range("D47").Value = 0.671
$ws.Range("E47").Value = 0.142
$ws.Range("F47").Value = 0.828

$ws.Range("C49").Value = 5.102
$ws.Range("D49").Value = 0.661
$ws.Range("E49").Value = 0.1
$ws.Range("F49").Value = 0.76

$ws.Range("C50").Value = 4.423
$ws.Range("D50").Value = 0.611
$ws.Range("E50").Value = 0.255
$ws.Range("F50").Value = 0.958

$ws.Range("C51").Value = 4.34
$ws.Range("D51").Value = 0.771
$ws.Range("E51").Value = 0.139
$ws.Range("F51").Value = 1.01

$ws.Range("C52").Value = 4.929
$ws.Range("D52").Value = 0.676
$ws.Range("E52").Value = 0.157
$ws.Range("F52").Value = 0.791

$ws.Range("C54").Value = 4.277
$ws.Range("D54").Value = 0.653
$ws.Range("E54").Value = 0.143
$ws.Range("F54").Value = 0.888

$ws.Range("C55").Value = 4.345
$ws.Range("D55").Value = 0.767
$ws.Range("E55").Value = 0.093
$ws.Range("F55").Value = 0.908

$ws.Range("C56").Value = 4.725
$ws.Range("D56").Value = 0.759
$ws.Range("E56").Value = 0.122
$ws.Range("F56").Value = 0.796

$ws.Range("C57").Value = 4.885
$ws.Range("D57").Value = 0.639
$ws.Range("E57").Value = 0.17
$ws.Range("F57").Value = 1.245

$ws.Range("C58").Value = 4.856
$ws.Range("D58").Value = 0.767
$ws.Range("E58").Value = 0.129
$ws.Range("F58").Value = 0.977

$ws.Range("C59").Value = 2.998
$ws.Range("D59").Value = 0.613
$ws.Range("E59").Value = 0.138
$ws.Range("F59").Value = 0.692

$ws.Range("C60").Value = 4.315
$ws.Range("D60").Value = 0.576
$ws.Range("E60").Value = 0.097
$ws.Range("F60").Value = 0.736

$ws.Range("C61").Value = 4.456538
$ws.Range("D61").Value = 0.680231
$ws.Range("E61").Value = 0.136692
$ws.Range("F61").Value = 0.832077

$ws.Range("C62").Value = 4.311245
$ws.Range("D62").Value = 0.611498
$ws.Range("E62").Value = 0.121188
$ws.Range("F62").Value = 0.91294

$ws.Range("C63").Value = 5.966618
$ws.Range("D63").Value = 0.779509
$ws.Range("E63").Value = 0.119217
$ws.Range("F63").Value = 1.447509

$ws.Range("C64").Value = 6.155478
$ws.Range("D64").Value = 0.821318
$ws.Range("E64").Value = 0.094909
$ws.Range("F64").Value = 1.420091

$ws.Range("C65").Value = 4.905357
$ws.Range("D65").Value = 0.718038
$ws.Range("E65").Value = 0.15
$ws.Range("F65").Value = 1.212115

$ws.Range("C66").Value = 5.471061
$ws.Range("D66").Value = 0.72944
$ws.Range("E66").Value = 0.119804
$ws.Range("F66").Value = 1.288201

$ws.Range("C67").Value = 4.118944
$ws.Range("D67").Value = 0.586
$ws.Range("E67").Value = 0.125588
$ws.Range("F67").Value = 0.919059

$ws.Range("C68").Value = 4.414132
$ws.Range("D68").Value = 0.620801
$ws.Range("E68").Value = 0.12812
$ws.Range("F68").Value = 0.82719

$ws.Range("C69").Value = 3.728514
$ws.Range("D69").Value = 0.548944
$ws.Range("E69").Value = 0.177276
$ws.Range("F69").Value = 0.710611

$ws.Range("C70").Value = 4.133767
$ws.Range("D70").Value = 0.528477
$ws.Range("E70").Value = 0.145189
$ws.Range("F70").Value = 0.629

$ws.Range("C71").Value = 4.456538
$ws.Range("D71").Value = 0.680231
$ws.Range("E71").Value = 0.136692
$ws.Range("F71").Value = 0.832077

$ws.Range("C72").Value = 4.0286
$ws.Range("D72").Value = 0.584322
$ws.Range("E72").Value = 0.193187
$ws.Range("F72").Value = 0.974056

$ws.Range("C73").Value = 4.110533
$ws.Range("D73").Value = 0.621571
$ws.Range("E73").Value = 0.095571
$ws.Range("F73").Value = 1.054786

$ws.Range("C74").Value = 5.0214
$ws.Range("D74").Value = 0.566
$ws.Range("E74").Value = 0.0688
$ws.Range("F74").Value = 1.09

$ws.Range("C75").Value = 4.4925
$ws.Range("D75").Value = 0.813
$ws.Range("E75").Value = 0.127
$ws.Range("F75").Value = 1.131

$ws.Range("C76").Value = 5.642333
$ws.Range("D76").Value = 0.904778
$ws.Range("E76").Value = 0.178889
$ws.Range("F76").Value = 1.343333

$ws.Range("C77").Value = 6.185667
$ws.Range("D77").Value = 0.805545
$ws.Range("E77").Value = 0.074
$ws.Range("F77").Value = 1.512091

$ws.Range("C78").Value = 6.582556
$ws.Range("D78").Value = 0.797
$ws.Range("E78").Value = 0.104
$ws.Range("F78").Value = 1.636926

$ws.Range("C79").Value = 6.702
$ws.Range("D79").Value = 0.807816
$ws.Range("E79").Value = 0.114158
$ws.Range("F79").Value = 1.632158

$ws.Range("C80").Value = 4.634825
$ws.Range("D80").Value = 0.502801
$ws.Range("E80").Value = 0.112848
$ws.Range("F80").Value = 1.009468

$ws.Range("C81").Value = 5.828368
$ws.Range("D81").Value = 0.754118
$ws.Range("E81").Value = 0.139118
$ws.Range("F81").Value = 1.443765

$ws.Range("C82").Value = 4.244869
$ws.Range("D82").Value = 0.631521
$ws.Range("E82").Value = 0.122724
$ws.Range("F82").Value = 0.895158

$ws.Range("C83").Value = 5.995484
$ws.Range("D83").Value = 0.78436
$ws.Range("E83").Value = 0.115416
$ws.Range("F83").Value = 1.448225

$ws.Range("C84").Value = 4.010974
$ws.Range("D84").Value = 0.580746
$ws.Range("E84").Value = 0.159192
$ws.Range("F84").Value = 0.731238

$ws.Range("C85").Value = 2.795667
$ws.Range("D85").Value = 0.2235
$ws.Range("E85").Value = 0.0705
$ws.Range("F85").Value = 0.741

$ws.Range("C86").Value = 4.358476
$ws.Range("D86").Value = 0.62465
$ws.Range("E86").Value = 0.10095
$ws.Range("F86").Value = 0.94205

$ws.Range("C87").Value = 5.02255
$ws.Range("D87").Value = 0.771474
$ws.Range("E87").Value = 0.143158
$ws.Range("F87").Value = 1.156158

$ws.Range("C88").Value = 5.165667
$ws.Range("D88").Value = 0.671833
$ws.Range("E88").Value = 0.0645
$ws.Range("F88").Value = 1.338

$ws.Range("C89").Value = 5.858829
$ws.Range("D89").Value = 0.759265
$ws.Range("E89").Value = 0.111706
$ws.Range("F89").Value = 1.413882

$ws.Range("C90").Value = 6.602902
$ws.Range("D90").Value = 0.82012
$ws.Range("E90").Value = 0.11648
$ws.Range("F90").Value = 1.60844

$ws.Range("C91").Value = 4.062362
$ws.Range("D91").Value = 0.598719
$ws.Range("E91").Value = 0.148313
$ws.Range("F91").Value = 0.765973

$ws.Range("C92").Value = 3.958125
$ws.Range("D92").Value = 0.607625
$ws.Range("E92").Value = 0.18425
$ws.Range("F92").Value = 0.815

$ws.Range("C93").Value = 4.793
$ws.Range("D93").Value = 0.5105
$ws.Range("E93").Value = 0.106
$ws.Range("F93").Value = 1.011

$ws.Range("C94").Value = 5.747375
$ws.Range("D94").Value = 0.743
$ws.Range("E94").Value = 0.164571
$ws.Range("F94").Value = 1.267143

$ws.Range("C95").Value = 3.816781
$ws.Range("D95").Value = 0.570026
$ws.Range("E95").Value = 0.133278
$ws.Range("F95").Value = 0.779205

$ws.Range("C96").Value = 5.371
$ws.Range("D96").Value = 0.7446
$ws.Range("E96").Value = 0.125667
$ws.Range("F96").Value = 1.344

$ws.Range("C97").Value = 4.210357
$ws.Range("D97").Value = 0.598528
$ws.Range("E97").Value = 0.141256
$ws.Range("F97").Value = 0.84319

$ws.Range("C98").Value = 4.535353
$ws.Range("D98").Value = 0.6005
$ws.Range("E98").Value = 0.159938
$ws.Range("F98").Value = 1.053875

# Fix mojibake in the PALOP footnote (A103): UTF-8 bytes that had been
# mis-decoded as Latin-1 (Ã­ -> í, Ã³ -> ó, etc.)
$ws.Range("A103").Value = 'Regional Economic Communities: CEN-SAD = "Community of Sahel-Saharan States"; COMESA = "Common Market for Eastern and Southern Africa"; EAC = "East African Community"; ECCAS = "Economic Community of Central African States"; ECOWAS = "Economic Community of West African States"; IGAD = "Intergovernmental Authority on Development"; SADC = "Southern African Development Community"; UMA = "Arab Maghreb Union"; PALOP = "Países Africanos de Língua Oficial Portuguesa"; ASEAN = "Association of Southeast Asian Nations"; MERCOSUR = "Mercado Común del Sur". EU27 = "European Union (27 members)". OECD = "Organisation for Economic Co-operation and Development".'
